# "Adicionei mais uma música"
#
# - Remove the proofErr (spell-check squiggle) wrapper around "Purpose".
# - Merge the run fragments "Hunting" / " High " / "and" / " " / "Low" into a
#   single run "Hunting High and Low" and drop the proofErr wrappers.
# - Add a new list item "Eyes Without a Face" after "Hunting High and Low",
#   using the same list-paragraph formatting as its siblings.
#
# Word only records w:proofErr spell-check markers as a side effect of
# editing text in place (Find/Replace, Range.Text=, etc. all leave the
# existing proofErr start/end markers sitting next to the run). Deleting a
# paragraph outright (range including its end-of-paragraph mark) drops its
# proofErr markers with it; InsertParagraphAfter() then manufactures a clean
# paragraph (inheriting formatting from its neighbour) with a single plain
# run and no proofErr markers at all. So: delete the two affected paragraphs
# wholesale, then rebuild all three (the two originals + the new song) from
# scratch after the "Forever Young" anchor paragraph.

$d = $word.ActiveDocument

$pPurpose = $d.Paragraphs.Item(2)
$pHunting = $d.Paragraphs.Item(3)

$toDelete = $d.Range($pPurpose.Range.Start, $pHunting.Range.End)
$toDelete.Delete()

$titles = @("Purpose", "Hunting High and Low", "Eyes Without a Face")
$anchor = $d.Paragraphs.Item(1)
foreach ($title in $titles) {
    $anchor.Range.InsertParagraphAfter()
    $anchor = $d.Paragraphs.Item($anchor.Index + 1)
    $anchor.Range.Text = $title
}
